$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '69.339.39'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  +0.02%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '3.419.42'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  +0.98%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  +0.06%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '580.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  -1.34%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '176.45'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  -2.40%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  +0.09%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'" + '3.412.77'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  +0.93%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'" + '0.591'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  -0.60%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'" + '0.197'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  +0.47%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  -1.04%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'" + '48.77'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  -0.39%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + '0.0000279'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  -1.95%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'" + '690.11'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  +1.10%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'" + '3.966.11'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  +0.90%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'" + '8.63'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  -0.02%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'" + '69.363.70'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  -0.01%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'" + '3.426.38'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  +1.35%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  +0.77%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'" + '17.65'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  -0.78%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'" + '11.37'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  -0.38%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  -0.67%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'" + '5.43'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  +0.61%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '16.91'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  -1.16%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'" + '100.53'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Value = "'" + '3.88'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  -1.57%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  -2.59%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  -1.04%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'" + '33.38'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  -2.92%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'" + '8.72'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  +0.21%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  -1.26%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'" + '569.52'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  +1.91%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  +0.45%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'" + '10.98'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  -2.05%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = "'" + '58.09'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  -0.01%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = "'" + '0.103'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  -3.06%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  -0.04%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'" + '3.580.46'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  -3.85%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  -2.01%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'" + '34.80'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  -0.71%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'" + '0.0₃0725'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  +2.44%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  -0.56%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'" + '2.65'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  -1.25%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'" + '0.331'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  -2.80%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'" + '0.0415'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  -0.69%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'" + '1.43'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  +2.49%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'" + '2.64'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  -1.09%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  -1.51%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'" + '0.999'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  -0.33%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'" + '131.89'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  -0.61%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'" + '2.66'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  +2.12%  '
$ws.Range('E51').Style = 'Normal'
